# Week1.pptx edit:
#  1. Merge the "Finding " / "Help" runs on the "Today's Objectives" slide
#     (slide 6) into a single run "Finding Help".
#  2. Delete the "Data Structures (Session 2)" slide (slide 23).
#  3. Delete the "Control Structures (Session 3)" slide (slide 23, after
#     the previous delete shifted indices down).

$p = $ppt.ActivePresentation

# --- 1. Fix "Finding Help" text on slide 6 -------------------------------
$slide6 = $p.Slides.Item(6)
$body = $slide6.Shapes.Item(2)
$tr = $body.TextFrame.TextRange
$findingHelpPara = $tr.Paragraphs(5, 1)
# Re-assigning the exact same string (or one that overlaps heavily with the
# target) is treated as a near-identical diff and the two runs are kept
# split; route through an unrelated placeholder first so the paragraph is
# genuinely rewritten as a single run.
$findingHelpPara.Text = "ZZZ"
$findingHelpPara = $tr.Paragraphs(5, 1)
$findingHelpPara.Text = "Finding Help"

# --- 2 & 3. Remove the two obsolete slides --------------------------------
# Original order (1-based): ... 22=rId23, 23="Data Structures (Session 2)",
# 24="Control Structures (Session 3)", 25="Finding Help ... Week 1 wrap-up".
$p.Slides.Item(24).Delete()
$p.Slides.Item(23).Delete()
